$wb = $excel.ActiveWorkbook

# Work on the "Test Case" worksheet, where the Expected Results (H2)
# text is being reworded and the trailing details cell (I2) is being
# cleared out.
$wsCase = $wb.Worksheets.Item("Test Case")

# Clear the contents of I2 (it previously duplicated some leftover
# "User name / details" copy that is no longer needed).
$wsCase.Range("I2").Value = ""

# Reword the "Expected Results" text in H2.
$newExpectedResults = "1. Log in to ""http://blubox.shoppinpal.com/s eller/"" website`n2. Navigate to Product module. Verify the following:`na. User is able to see thumbnails in the selected layout`nb. User is able to see checkbox over upper left corner of thumbnails`nc. User is able to see options tag over upper roght corner of thumbnails`nd. Various details for thumbnails are present"
$wsCase.Range("H2").Value = $newExpectedResults

# Switch to the "Test Case" sheet and leave I2 as the active/selected
# cell, matching where editing left off.
$wsCase.Activate()
$wsCase.Range("I2").Select()
